$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number and week-covering dates) ---
$ws.Range("C8").Value = "Volume 30   Number  28"
$ws.Range("C9").Value = "Report Covering the Week  7/10/2023  Through  7/16/2023"

# --- Cells whose number format/style changes from placeholder text to numeric ---
# style 15 -> "#,##0"   style 16 -> '#,##0.0;"-"#,##0.0'
$ws.Range("D15").Value = 1
$ws.Range("D15").NumberFormat = "#,##0"
$ws.Range("E15").Value = -100
$ws.Range("E15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("C23").Value = 3
$ws.Range("C23").NumberFormat = "#,##0"
$ws.Range("D27").Value = 1
$ws.Range("D27").NumberFormat = "#,##0"
$ws.Range("E27").Value = 0
$ws.Range("E27").NumberFormat = "#,##0.0;""-""#,##0.0"

# --- Plain numeric value updates (style/format unchanged) ---
$values = @{
    "G14" = 1
    "M14" = -85.714285714285
    "G15" = 2
    "J15" = 21
    "K15" = -38.095238095238
    "N15" = -69.047619047619
    "C16" = 12
    "D16" = 11
    "E16" = 9.090909090909
    "F16" = 25
    "G16" = 36
    "H16" = -30.555555555555
    "I16" = 184
    "J16" = 211
    "K16" = -12.796208530805
    "L16" = 24.324324324324
    "M16" = -28.404669260700
    "N16" = -78.923253150057
    "C17" = 8
    "D17" = 14
    "E17" = -42.857142857142
    "F17" = 54
    "G17" = 51
    "H17" = 5.882352941176
    "I17" = 366
    "J17" = 320
    "K17" = 14.375
    "L17" = 47.580645161290
    "M17" = 109.142857142857
    "N17" = -20.607375271149
    "C18" = 9
    "D18" = 5
    "E18" = 80
    "F18" = 19
    "G18" = 22
    "H18" = -13.636363636363
    "I18" = 95
    "J18" = 102
    "K18" = -6.862745098039
    "L18" = 25
    "M18" = -28.030303030303
    "N18" = -84.627831715210
    "C19" = 17
    "D19" = 7
    "E19" = 142.857142857143
    "F19" = 45
    "G19" = 40
    "H19" = 12.5
    "I19" = 293
    "J19" = 303
    "K19" = -3.300330033003
    "L19" = 38.207547169811
    "M19" = 50.256410256410
    "N19" = -53.044871794871
    "C20" = 6
    "E20" = 100
    "F20" = 19
    "G20" = 12
    "H20" = 58.333333333333
    "I20" = 129
    "J20" = 130
    "K20" = -0.769230769230
    "L20" = 118.64406779661
    "M20" = 27.722772277227
    "N20" = -83.875
    "C21" = 52
    "D21" = 41
    "E21" = 26.829268292682
    "F21" = 162
    "G21" = 164
    "H21" = -1.219512195121
    "I21" = 1081
    "J21" = 1093
    "K21" = -1.097895699908
    "L21" = 42.800528401585
    "M21" = 23.401826484018
    "N21" = -68.502331002331
    "L22" = 233.333333333333
    "F23" = 6
    "H23" = 500
    "I23" = 23
    "K23" = -11.538461538461
    "L23" = 21.052631578947
    "M23" = 4.545454545454
    "C24" = 32
    "E24" = -15.789473684210
    "F24" = 132
    "G24" = 129
    "H24" = 2.325581395348
    "I24" = 853
    "J24" = 811
    "K24" = 5.178791615289
    "L24" = 52.321428571428
    "M24" = 32.248062015503
    "C25" = 25
    "D25" = 21
    "E25" = 19.047619047619
    "F25" = 87
    "G25" = 66
    "H25" = 31.818181818181
    "I25" = 536
    "J25" = 393
    "K25" = 36.386768447837
    "L25" = 71.246006389776
    "M25" = 22.654462242562
    "G26" = 3
    "J26" = 30
    "K26" = -30
    "C27" = 1
    "F27" = 7
    "G27" = 2
    "H27" = 250
    "I27" = 47
    "J27" = 42
    "K27" = 11.904761904761
    "L27" = 56.666666666666
    "L28" = -60
    "N28" = -92.405063291139
    "L29" = -50
    "N29" = -91.891891891891
}
foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
